$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# Collapse "mieten. Die Kunden werden über ein Transportweg (Luft/Land/Wasser)
# befördert. Ein Kunde kann je Auftrag ein Tier buchen. " down to
# "mieten. Ein Kunde kann je Auftrag ein Tier buchen. "
$rng1 = $d.Content
$old1 = "mieten. Die Kunden werden über ein Transportweg (Luft/Land/Wasser) befördert. Ein Kunde kann je Auftrag ein Tier buchen. "
$new1 = "mieten. Ein Kunde kann je Auftrag ein Tier buchen. "
$rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- Change 2 ---------------------------------------------------------
# Merge the five separate runs that make up
# "Bei Zahlungseingang wird der Bestellstatus manuell auf „bezahlt“ gesetzt."
# into a single run (same visible text, replaced in place so the engine
# coalesces the runs).
$rng2 = $d.Content
$old2 = "Bei Zahlungseingang wird der Bestellstatus manuell auf „bezahlt“ gesetzt."
$rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2)

# --- Change 3 ---------------------------------------------------------
# Add one more empty paragraph right after that paragraph (before the
# existing run of empty paragraphs at the end of the document).
$rng3 = $d.Content
$rng3.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.Collapse(0)
$rng3.Text = "`r"
